# Apply the "break out stock.yaml completed" edit to the weekly PIIND.NS sheet.
#
# Summary of changes:
#  - Q56:  1 -> 0
#  - Q60:  2 -> 0
#  - O679: 0 -> 2
#  - R681: (blank) -> 0
#  - R682: (blank) -> 0
#  - Append 17 new weekly OHLCV rows (683-699), matching the existing table
#    layout (A:Datetime-serial .. R:backup). New R cells are left blank,
#    matching the pre-existing pattern for the most recent rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- targeted single-cell fixes -------------------------------------------
$ws.Range("Q56").Value = 0
$ws.Range("Q60").Value = 0
$ws.Range("O679").Value = 2
$ws.Range("R681").Value = 0
$ws.Range("R682").Value = 0

# --- append new rows 683-699 ----------------------------------------------
# Columns: A Datetime, B Open, C High, D Low, E Close, F Adj Close, G Volume,
#          H Year, I Month, J Day, K Hour, L Minute, M Second, N Week,
#          O isPivot, P two_line_structure, Q detect_structure
$newRows = @(
    @(683,45474,3799.5,3845.85009765625,3728,3788.5,3780.614013671875,1312657,2024,7,1,0,0,0,27,0,0,0),
    @(684,45481,3765.14990234375,3922.199951171875,3730.10009765625,3901.300048828125,3893.17919921875,1849392,2024,7,8,0,0,0,28,0,0,0),
    @(685,45488,3925.699951171875,3928.300048828125,3812,3835.300048828125,3827.316650390625,1062468,2024,7,15,0,0,0,29,0,0,1),
    @(686,45495,3849.949951171875,4214.5,3777,4207,4198.24267578125,2483889,2024,7,22,0,0,0,30,0,0,0),
    @(687,45502,4208.25,4467.0498046875,4208.25,4278.89990234375,4269.9931640625,2153698,2024,7,29,0,0,0,31,0,0,0),
    @(688,45509,4220,4600,4185.14990234375,4487.7001953125,4478.35888671875,2753757,2024,8,5,0,0,0,32,0,0,0),
    @(689,45516,4484.85009765625,4517.85009765625,4294,4436.5,4427.26513671875,1184436,2024,8,12,0,0,0,33,0,0,0),
    @(690,45523,4450,4497.10009765625,4257,4420.5,4411.29833984375,1768015,2024,8,19,0,0,0,34,0,0,0),
    @(691,45530,4449.25,4526.9501953125,4383.7001953125,4495.4501953125,4495.4501953125,1379162,2024,8,26,0,0,0,35,0,0,0),
    @(692,45537,4524.89990234375,4660.9501953125,4476.60009765625,4613.25,4613.25,1615043,2024,9,2,0,0,0,36,0,0,0),
    @(693,45544,4613.25,4748.89990234375,4565.7001953125,4660.25,4660.25,1412314,2024,9,9,0,0,0,37,0,0,0),
    @(694,45551,4675,4804.0498046875,4598.2001953125,4716,4716,1437849,2024,9,16,0,0,0,38,1,0,0),
    @(695,45558,4746.60009765625,4765.9501953125,4522.64990234375,4652.9501953125,4652.9501953125,963857,2024,9,23,0,0,0,39,0,0,0),
    @(696,45565,4652.9501953125,4759,4570.5498046875,4587.25,4587.25,707949,2024,9,30,0,0,0,40,0,0,0),
    @(697,45572,4590,4657.60009765625,4472.7998046875,4523.60009765625,4523.60009765625,578128,2024,10,7,0,0,0,41,0,0,0),
    @(698,45579,4543.4501953125,4618.89990234375,4407.0498046875,4560.7001953125,4560.7001953125,806309,2024,10,14,0,0,0,42,0,0,0),
    @(699,45586,4560.89990234375,4635.10009765625,4259,4327.64990234375,4327.64990234375,833902,2024,10,21,0,0,0,43,0,0,0)
)

foreach ($row in $newRows) {
    $r = $row[0]
    for ($c = 1; $c -le 17; $c++) {
        $ws.Cells.Item($r, $c).Value = $row[$c]
    }
    # Column A carries the same datetime number format as the rest of the
    # table (style index 2 in the original workbook).
    $ws.Cells.Item($r, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    # Column R (backup) is left blank for these rows, same as rows 681/682
    # were before this edit.
}

Write-Host "Applied PIIND.NS weekly sheet update: fixed Q56/Q60/O679/R681/R682 and appended rows 683-699"
